$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "275.11"
Set-TextValue "E2" "0.40%"
Set-TextValue "E3" "1.96%"
Set-TextValue "D4" "4.847"
Set-TextValue "E4" "-0.54%"
Set-TextValue "D5" "0.06389"
Set-TextValue "E5" "1.01%"
Set-TextValue "E6" "0.93%"
Set-TextValue "E7" "-2.80%"
Set-TextValue "D8" "0.8776"
Set-TextValue "E8" "0.68%"
Set-TextValue "E9" "3.74%"
Set-TextValue "E10" "0.30%"
Set-TextValue "D11" "0.07569"
Set-TextValue "E11" "2.54%"
Set-TextValue "D12" "0.02972"
Set-TextValue "E12" "-1.37%"
Set-TextValue "D13" "0.08976"
Set-TextValue "E13" "-0.74%"
Set-TextValue "D14" "0.001566"
Set-TextValue "E14" "-0.28%"
Set-TextValue "D15" "0.0006397"
Set-TextValue "E15" "1.69%"
Set-TextValue "D16" "0.006184"
Set-TextValue "E16" "4.20%"
Set-TextValue "D17" "3.472"
Set-TextValue "E17" "0.57%"
Set-TextValue "D18" "3.309"
Set-TextValue "E18" "-0.64%"
Set-TextValue "D19" "2.251"
Set-TextValue "E19" "-1.43%"
Set-TextValue "E20" "-0.95%"
Set-TextValue "D21" "0.1350"
Set-TextValue "E21" "1.90%"
Set-TextValue "D22" "3.909"
Set-TextValue "E22" "0.09%"
Set-TextValue "D23" "0.04403"
Set-TextValue "E23" "1.18%"
Set-TextValue "D25" "0.001179"
Set-TextValue "E25" "0.38%"
Set-TextValue "D26" "0.003851"
Set-TextValue "E26" "-9.86%"
Set-TextValue "E27" "0.13%"
Set-TextValue "E28" "14.68%"
Set-TextValue "D40" "0.04131"
Set-TextValue "E40" "2.52%"
Set-TextValue "D41" "0.006818"
Set-TextValue "E41" "1.37%"
Set-TextValue "D42" "0.1173"
Set-TextValue "E42" "0.53%"
Set-TextValue "D43" "0.002151"
Set-TextValue "E43" "2.51%"
Set-TextValue "D44" "0.01148"
Set-TextValue "E44" "-8.44%"
Set-TextValue "D45" "0.00005202"
Set-TextValue "E45" "-1.80%"
Set-TextValue "D46" "1.680"
Set-TextValue "E46" "-28.80%"
Set-TextValue "D47" "0.01855"
Set-TextValue "E47" "-7.22%"
